# Updated cryptos list on Mon Feb 12 06:19:14 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, [string]$Text)
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.Style = "Normal"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "48.087.42"
$ws.Range("E2").Value = "  -0.41%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.495.38"
$ws.Range("E3").Value = "  -1.18%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.10%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "318.78"
$ws.Range("E5").Value = "  -1.40%  "

# Row 6 - Solana
Set-TextValue $ws.Range("D6") "105.56"
$ws.Range("E6").Value = "  -3.28%  "

# Row 7 - XRP
Set-TextValue $ws.Range("D7") "0.519"
$ws.Range("E7").Value = "  -1.55%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.05%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  -4.66%  "

# Row 10 - Avalanche
Set-TextValue $ws.Range("D10") "38.76"
$ws.Range("E10").Value = "  -4.11%  "

# Row 11 - Chainlink
Set-TextValue $ws.Range("D11") "20.19"
$ws.Range("E11").Value = "  -0.17%  "

# Row 12 - Dogecoin
$ws.Range("E12").Value = "  -2.27%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  -0.41%  "

# Row 14 - Polkadot
Set-TextValue $ws.Range("D14") "7.07"
$ws.Range("E14").Value = "  -3.07%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "2.885.92"
$ws.Range("E15").Value = "  -1.10%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "2.495.73"
$ws.Range("E16").Value = "  -1.43%  "

# Row 17 - Polygon
$ws.Range("E17").Value = "  -3.79%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "47.922.80"
$ws.Range("E18").Value = "  -0.51%  "

# Row 19 - ImmutableX
Set-TextValue $ws.Range("D19") "2.98"
$ws.Range("E19").Value = "  +9.16%  "

# Row 20 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D20") "12.92"
$ws.Range("E20").Value = "  -2.59%  "

# Row 21 - Uniswap
$ws.Range("E21").Value = "  -0.49%  "

# Row 22 - ShibaInu
$ws.Range("D22").Value = "0.0₃0929"
$ws.Range("E22").Value = "  -2.01%  "

# Row 23 - Litecoin
Set-TextValue $ws.Range("D23") "71.01"
$ws.Range("E23").Value = "  -2.05%  "

# Row 24 - BitcoinCash
Set-TextValue $ws.Range("D24") "270.82"
$ws.Range("E24").Value = "  -0.67%  "

# Row 25 - PancakeSwap
$ws.Range("E25").Value = "  -2.70%  "

# Row 26 - Dai
$ws.Range("E26").Value = "  +0.04%  "

# Row 27 - EthereumClassic
Set-TextValue $ws.Range("D27") "25.74"
$ws.Range("E27").Value = "  -2.00%  "

# Row 28 - Toncoin
Set-TextValue $ws.Range("D28") "2.29"
$ws.Range("E28").Value = "  -0.73%  "

# Row 29 - Cosmos
Set-TextValue $ws.Range("D29") "9.71"
$ws.Range("E29").Value = "  -4.59%  "

# Row 30 - Kaspa
Set-TextValue $ws.Range("D30") "0.138"
$ws.Range("E30").Value = "  -5.74%  "

# Row 31 - InjectiveProtocol
Set-TextValue $ws.Range("D31") "34.45"
$ws.Range("E31").Value = "  -2.21%  "

# Row 32 - OKB
Set-TextValue $ws.Range("D32") "49.37"
$ws.Range("E32").Value = "  -0.83%  "

# Row 33 - FirstDigitalUSD
$ws.Range("E33").Value = "  -0.09%  "

# Row 34 - Celestia(->Filecoin)
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D34") "5.28"
$ws.Range("E34").Value = "  -2.39%  "

# Row 35 - Filecoin(->Celestia)
$ws.Range("B35").Value = "Celestia"
$ws.Range("C35").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextValue $ws.Range("D35") "18.94"
$ws.Range("E35").Value = "  -5.37%  "

# Row 36 - Hedera
Set-TextValue $ws.Range("D36") "0.0772"
$ws.Range("E36").Value = "  -2.65%  "

# Row 37 - ARBITRUM
$ws.Range("E37").Value = "  -3.11%  "

# Row 38 - RenderToken
$ws.Range("E38").Value = "  -3.67%  "

# Row 39 - LidoDAOToken
$ws.Range("E39").Value = "  -4.66%  "

# Row 40 - Monero
Set-TextValue $ws.Range("D40") "121.77"
$ws.Range("E40").Value = "  +2.69%  "

# Row 41 - WEMIXToken(->Stellar)
$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws.Range("D41") "0.110"
$ws.Range("E41").Value = "  -2.30%  "

# Row 42 - Stellar(->EnergySwap)
$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D42") "22.22"
$ws.Range("E42").Value = "  -0.49%  "

# Row 43 - EnergySwap(->WEMIXToken)
$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Range("D43") "2.22"
$ws.Range("E43").Value = "  +1.01%  "

# Row 44 - VeChain
Set-TextValue $ws.Range("D44") "0.0301"
$ws.Range("E44").Value = "  +0.48%  "

# Row 45 - Maker
$ws.Range("D45").Value = "1.997.40"
$ws.Range("E45").Value = "  -0.55%  "

# Row 46 - NEARProtocol
Set-TextValue $ws.Range("D46") "3.12"
$ws.Range("E46").Value = "  -1.71%  "

# Row 47 - Stacks
Set-TextValue $ws.Range("D47") "1.87"
$ws.Range("E47").Value = "  -1.77%  "

# Row 49 - FraxShare
Set-TextValue $ws.Range("D49") "8.93"
$ws.Range("E49").Value = "  -1.94%  "

# Row 50 - THORChain
Set-TextValue $ws.Range("D50") "5.15"
$ws.Range("E50").Value = "  -1.76%  "

# Row 51 - BitcoinSV
Set-TextValue $ws.Range("D51") "78.97"
$ws.Range("E51").Value = "  -2.00%  "
